$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t = $d.Tables.Item(1)
$cell = $t.Cell(7, 7)
$cell.Range.Text = "19 (63)`r15 (37)`r9 (75)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(7, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(7, 8)
$cell.Range.Text = "0* (1)`r5 (11)`r0"
$t = $d.Tables.Item(1)
$cell = $t.Cell(7, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(8, 7)
$cell.Range.Text = "11 (42)`r9 (18)`r3 (25)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(8, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(8, 8)
$cell.Range.Text = "0`r1 (3)`r0"
$t = $d.Tables.Item(1)
$cell = $t.Cell(8, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(9, 7)
$cell.Range.Text = "0`r0`r0"
$t = $d.Tables.Item(1)
$cell = $t.Cell(9, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(9, 8)
$cell.Range.Text = "0`r0`r0"
$t = $d.Tables.Item(1)
$cell = $t.Cell(9, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(10, 7)
$cell.Range.Text = "0`r0`r2 (3)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(10, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(10, 8)
$cell.Range.Text = "0`r0`r0"
$t = $d.Tables.Item(1)
$cell = $t.Cell(10, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(11, 7)
$cell.Range.Text = "0`r1 (2)`r0"
$t = $d.Tables.Item(1)
$cell = $t.Cell(11, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(11, 8)
$cell.Range.Text = "0`r0`r0"
$t = $d.Tables.Item(1)
$cell = $t.Cell(11, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(12, 7)
$cell.Range.Text = "2 (4)`r1 (2)`r0"
$t = $d.Tables.Item(1)
$cell = $t.Cell(12, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(12, 8)
$cell.Range.Text = "0`r0`r0"
$t = $d.Tables.Item(1)
$cell = $t.Cell(12, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(13, 7)
$cell.Range.Text = "85 (155)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(13, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(13, 8)
$cell.Range.Text = "4 (6)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(13, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(14, 7)
$cell.Range.Text = "59 (88)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(14, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(14, 8)
$cell.Range.Text = "9 (13)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(14, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(15, 7)
$cell.Range.Text = "83 (120)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(15, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(15, 8)
$cell.Range.Text = "0"
$t = $d.Tables.Item(1)
$cell = $t.Cell(15, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(16, 7)
$cell.Range.Text = "58 (112)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(16, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(16, 8)
$cell.Range.Text = "3 (6)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(16, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(17, 7)
$cell.Range.Text = "94 (79)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(17, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(17, 8)
$cell.Range.Text = "0"
$t = $d.Tables.Item(1)
$cell = $t.Cell(17, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(18, 7)
$cell.Range.Text = "91 (103)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(18, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(18, 8)
$cell.Range.Text = "1 (1)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(18, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(19, 7)
$cell.Range.Text = "68 (90)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(19, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(19, 8)
$cell.Range.Text = "8 (10)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(19, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(20, 7)
$cell.Range.Text = "61 (82)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(20, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(20, 8)
$cell.Range.Text = "4 (5)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(20, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(21, 7)
$cell.Range.Text = "73 (82)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(21, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(21, 8)
$cell.Range.Text = "2 (2)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(21, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(22, 7)
$cell.Range.Text = "55 (72)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(22, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(22, 8)
$cell.Range.Text = "5 (5)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(22, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(23, 7)
$cell.Range.Text = "52 (81)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(23, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(23, 8)
$cell.Range.Text = "1 (1)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(23, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(24, 7)
$cell.Range.Text = "47 (74)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(24, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(24, 8)
$cell.Range.Text = "4 (7)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(24, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(25, 7)
$cell.Range.Text = "92 (94)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(25, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(25, 8)
$cell.Range.Text = "3 (3)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(25, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(26, 7)
$cell.Range.Text = "77 (93)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(26, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(26, 8)
$cell.Range.Text = "4 (5)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(26, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(27, 7)
$cell.Range.Text = "77 (101)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(27, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(27, 8)
$cell.Range.Text = "1 (1)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(27, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(28, 7)
$cell.Range.Text = "68 (79)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(28, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(28, 8)
$cell.Range.Text = "4 (5)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(28, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(29, 7)
$cell.Range.Text = "61 (54)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(29, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(29, 8)
$cell.Range.Text = "3 (2)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(29, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(30, 7)
$cell.Range.Text = "72 (63)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(30, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(30, 8)
$cell.Range.Text = "4 (4)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(30, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(31, 7)
$cell.Range.Text = "48 (53)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(31, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(31, 8)
$cell.Range.Text = "2 (3)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(31, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(32, 7)
$cell.Range.Text = "72 (86)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(32, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(32, 8)
$cell.Range.Text = "5 (5)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(32, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(33, 7)
$cell.Range.Text = "64 (80)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(33, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(33, 8)
$cell.Range.Text = "5 (6)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(33, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(34, 7)
$cell.Range.Text = "67 (84)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(34, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(34, 8)
$cell.Range.Text = "2 (2)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(34, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(35, 7)
$cell.Range.Text = "7 (39)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(35, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(35, 8)
$cell.Range.Text = "0"
$t = $d.Tables.Item(1)
$cell = $t.Cell(35, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(36, 7)
$cell.Range.Text = "56 (55)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(36, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(36, 8)
$cell.Range.Text = "7 (7)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(36, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(37, 7)
$cell.Range.Text = "84 (87)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(37, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(37, 8)
$cell.Range.Text = "8 (9)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(37, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(38, 7)
$cell.Range.Text = "88 (84)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(38, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(38, 8)
$cell.Range.Text = "5 (4)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(38, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(39, 7)
$cell.Range.Text = "88 (86)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(39, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(39, 8)
$cell.Range.Text = "9 (9)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(39, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(40, 7)
$cell.Range.Text = "83 (69)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(40, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(40, 8)
$cell.Range.Text = "6 (5)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(40, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(41, 7)
$cell.Range.Text = "81 (76)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(41, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(41, 8)
$cell.Range.Text = "4 (4)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(41, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(42, 7)
$cell.Range.Text = "25 (36)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(42, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(42, 8)
$cell.Range.Text = "1 (1)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(42, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(43, 3)
$cell.Range.Text = "0"
$t = $d.Tables.Item(1)
$cell = $t.Cell(43, 3)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(43, 4)
$cell.Range.Text = "11 (15)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(43, 4)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(43, 5)
$cell.Range.Text = "14 (21)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(43, 5)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(43, 6)
$cell.Range.Text = "4 (6)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(43, 6)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(43, 7)
$cell.Range.Text = "31 (43)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(43, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(43, 8)
$cell.Range.Text = "15 (21)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(43, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(43, 9)
$cell.Range.Text = "72"
$t = $d.Tables.Item(1)
$cell = $t.Cell(43, 9)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(44, 7)
$cell.Range.Text = "61 (79)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(44, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(44, 8)
$cell.Range.Text = "8 (11)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(44, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(45, 7)
$cell.Range.Text = "85 (88)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(45, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(45, 8)
$cell.Range.Text = "7 (8)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(45, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(46, 7)
$cell.Range.Text = "52 (65)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(46, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(46, 8)
$cell.Range.Text = "9 (12)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(46, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(47, 7)
$cell.Range.Text = "76 (100)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(47, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(47, 8)
$cell.Range.Text = "6 (9)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(47, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(48, 7)
$cell.Range.Text = "81 (94)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(48, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(48, 8)
$cell.Range.Text = "7 (9)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(48, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(49, 7)
$cell.Range.Text = "73 (94)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(49, 7)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

$t = $d.Tables.Item(1)
$cell = $t.Cell(49, 8)
$cell.Range.Text = "3 (4)"
$t = $d.Tables.Item(1)
$cell = $t.Cell(49, 8)
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 10
$cell.Range.Font.SizeBi = 10

